# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The table of worker mora periods (rows 16-26, columns C..G) is
# refreshed with new content:
#   - JHAN KEVIN MONTALVO FONTALVO (CC 1047409809), period 1711 now
#     appears first (row 16).
#   - YEINER DE JESUS LOPEZ SANCHEZ (CC 1140841233) keeps all of his
#     periods (rows 17-26) but now listed in ascending order
#     (2103 -> 2112) instead of descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "1047409809"; Name = "JHAN KEVIN MONTALVO FONTALVO"; Period = "1711"; Mora = 984;   Salario = 781242 },
    @{ Row = 17; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2103"; Mora = 38050; Salario = 951231 },
    @{ Row = 18; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2104"; Mora = 38050; Salario = 951231 },
    @{ Row = 19; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2105"; Mora = 38050; Salario = 951231 },
    @{ Row = 20; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2106"; Mora = 38050; Salario = 951231 },
    @{ Row = 21; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2107"; Mora = 38050; Salario = 951231 },
    @{ Row = 22; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2108"; Mora = 38050; Salario = 951231 },
    @{ Row = 23; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2109"; Mora = 38050; Salario = 951231 },
    @{ Row = 24; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2110"; Mora = 38050; Salario = 951231 },
    @{ Row = 25; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2111"; Mora = 38050; Salario = 951231 },
    @{ Row = 26; Doc = "1140841233"; Name = "YEINER DE JESUS LOPEZ SANCHEZ"; Period = "2112"; Mora = 20293; Salario = 951231 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
